$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns keep their text formatting
# so values such as "326.60" or "-1.09%" are stored as literal text,
# matching the original inline-string cell contents rather than being
# auto-converted to numbers/percentages by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "326.60"
$ws.Range("E2").Value = "-1.09%"
$ws.Range("D3").Value = "39.63"
$ws.Range("E3").Value = "-1.11%"
$ws.Range("D4").Value = "5.701"
$ws.Range("E4").Value = "5.77%"
$ws.Range("D5").Value = "0.08044"
$ws.Range("E5").Value = "-1.10%"
$ws.Range("D6").Value = "2.007"
$ws.Range("E6").Value = "4.13%"
$ws.Range("D7").Value = "8.636"
$ws.Range("E7").Value = "-0.22%"
$ws.Range("D8").Value = "4.492"
$ws.Range("E8").Value = "-0.77%"
$ws.Range("D9").Value = "2.958"
$ws.Range("E9").Value = "-0.65%"
$ws.Range("D10").Value = "0.9232"
$ws.Range("E10").Value = "-2.19%"
$ws.Range("D11").Value = "0.1258"
$ws.Range("E11").Value = "-7.72%"
$ws.Range("D12").Value = "0.1967"
$ws.Range("E12").Value = "-0.64%"
$ws.Range("D13").Value = "8.747"
$ws.Range("E13").Value = "21.58%"
$ws.Range("D14").Value = "0.09217"
$ws.Range("E14").Value = "-1.28%"
$ws.Range("D15").Value = "0.03561"
$ws.Range("E15").Value = "0.21%"
$ws.Range("E16").Value = "9.51%"
$ws.Range("D17").Value = "0.001293"
$ws.Range("E17").Value = "-2.20%"
$ws.Range("D18").Value = "0.006311"
$ws.Range("E18").Value = "-1.49%"
$ws.Range("E19").Value = "0.04%"
$ws.Range("D20").Value = "0.3479"
$ws.Range("E20").Value = "-1.29%"
$ws.Range("E21").Value = "1.33%"
$ws.Range("D22").Value = "0.2699"
$ws.Range("E22").Value = "5.54%"
$ws.Range("E23").Value = "-0.68%"
$ws.Range("E24").Value = "2.98%"
$ws.Range("D25").Value = "0.004619"
$ws.Range("E25").Value = "7.52%"
$ws.Range("E26").Value = "-0.87%"
$ws.Range("D39").Value = "0.02500"
$ws.Range("E39").Value = "-0.59%"
$ws.Range("D40").Value = "0.05339"
$ws.Range("E40").Value = "1.96%"
$ws.Range("D41").Value = "0.007453"
$ws.Range("E41").Value = "-2.18%"
$ws.Range("D42").Value = "0.009899"
$ws.Range("E42").Value = "8.42%"
$ws.Range("E43").Value = "-1.49%"
$ws.Range("D44").Value = "0.002114"
$ws.Range("E44").Value = "-2.59%"
$ws.Range("D45").Value = "0.01088"
$ws.Range("E45").Value = "-0.24%"
$ws.Range("D46").Value = "0.00006673"
$ws.Range("E46").Value = "1.27%"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").Value = "-0.15%"
$ws.Range("D48").Value = "0.003037"
$ws.Range("E48").Value = "-9.19%"
$ws.Range("E49").Value = "-5.18%"
$ws.Range("D50").Value = "0.00002097"
$ws.Range("E50").Value = "-0.15%"
$ws.Range("D51").Value = "0.0001997"
$ws.Range("E51").Value = "-0.15%"
